$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 8: rename the "driving way" label to "driving distance"
$ws.Range("A8").Value = "single driving distance in meters"

# New row 9: number of iterations
$ws.Range("A9").Value = "number of iterations"
$ws.Range("A3").Copy()
$ws.Range("A9").PasteSpecial(-4122)

# Units column: spell out the abbreviations
$ws.Range("C6").Value = "seconds"
$ws.Range("C7").Value = "seconds"
$ws.Range("C8").Value = "meters"

# New B9 value (with border, like the other B cells)
$ws.Range("B9").Value = 50
$ws.Range("B3").Copy()
$ws.Range("B9").PasteSpecial(-4122)

# Widen column C to fit the new, longer unit text
$ws.Columns.Item(3).ColumnWidth = 14

# Update selection to reflect the new active cell location
$ws.Range("B11").Select()
